$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply formatting to the new cells by copying formats from existing,
#     already-formatted cells (this reuses existing style indices instead
#     of minting new ones). ---
$ws.Range("A2").Copy()
$ws.Range("A223:A232").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B223:B232").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C223:C232").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("D223:D232").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("E223:E232").PasteSpecial(-4122)

$ws.Range("K222").Copy()
$ws.Range("K223:K227").PasteSpecial(-4122)

$ws.Range("G2").Copy()
$ws.Range("G228:G232").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Row 223-227: Lake Florida / Austen, transect 11-15, Attached to "0 mussels" ---
$ws.Range("A223").Value = 43328.50091184028
$ws.Range("B223").Value = 43327.0
$ws.Range("C223").Value = "Lake Florida"
$ws.Range("D223").Value = "Austen"
$ws.Range("E223").Value = 11.0
$ws.Range("K223").Value = "0 mussels"

$ws.Range("A224").Value = 43328.50115063657
$ws.Range("B224").Value = 43327.0
$ws.Range("C224").Value = "Lake Florida"
$ws.Range("D224").Value = "Austen"
$ws.Range("E224").Value = 12.0
$ws.Range("K224").Value = "0 mussels"

$ws.Range("A225").Value = 43328.50136325232
$ws.Range("B225").Value = 43327.0
$ws.Range("C225").Value = "Lake Florida"
$ws.Range("D225").Value = "Austen"
$ws.Range("E225").Value = 13.0
$ws.Range("K225").Value = "0 mussels"

$ws.Range("A226").Value = 43328.50157587963
$ws.Range("B226").Value = 43327.0
$ws.Range("C226").Value = "Lake Florida"
$ws.Range("D226").Value = "Austen"
$ws.Range("E226").Value = 14.0
$ws.Range("K226").Value = "0 mussels"

$ws.Range("A227").Value = 43328.50180258101
$ws.Range("B227").Value = 43327.0
$ws.Range("C227").Value = "Lake Florida"
$ws.Range("D227").Value = "Austen"
$ws.Range("E227").Value = 15.0
$ws.Range("K227").Value = "0 mussels"

# --- Row 228-232: Lake Florida / Aislyn, transect 11-15, 0 mussels in cluster ---
$ws.Range("A228").Value = 43328.5456141551
$ws.Range("B228").Value = 43327.0
$ws.Range("C228").Value = "Lake Florida"
$ws.Range("D228").Value = "Aislyn"
$ws.Range("E228").Value = 11.0
$ws.Range("G228").Value = 0.0

$ws.Range("A229").Value = 43328.545763449074
$ws.Range("B229").Value = 43327.0
$ws.Range("C229").Value = "Lake Florida"
$ws.Range("D229").Value = "Aislyn"
$ws.Range("E229").Value = 12.0
$ws.Range("G229").Value = 0.0

$ws.Range("A230").Value = 43328.54590118055
$ws.Range("B230").Value = 43327.0
$ws.Range("C230").Value = "Lake Florida"
$ws.Range("D230").Value = "Aislyn"
$ws.Range("E230").Value = 13.0
$ws.Range("G230").Value = 0.0

$ws.Range("A231").Value = 43328.54602976852
$ws.Range("B231").Value = 43327.0
$ws.Range("C231").Value = "Lake Florida"
$ws.Range("D231").Value = "Aislyn"
$ws.Range("E231").Value = 14.0
$ws.Range("G231").Value = 0.0

$ws.Range("A232").Value = 43328.54617332176
$ws.Range("B232").Value = 43327.0
$ws.Range("C232").Value = "Lake Florida"
$ws.Range("D232").Value = "Aislyn"
$ws.Range("E232").Value = 15.0
$ws.Range("G232").Value = 0.0

Write-Output "Added rows 223-232 (Lake Florida sample)"
